$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$A2 = @'
Pipeline(steps=[('scaler', None),
                ('selector',
                 SelectFromModel(estimator=ExtraTreesClassifier(random_state=42))),
                ('model',
                 BaggingClassifier(estimator=SVC(C=0.001, kernel='linear',
                                                 random_state=42),
                                   random_state=42))])
'@
$ws.Range("A2").Value = $A2

$ws.Range("B2").Value = 0.6761904761904762
$C2 = @'
{'selector': SelectFromModel(estimator=ExtraTreesClassifier(random_state=42)), 'scaler': None, 'model__n_estimators': 10, 'model__estimator__kernel': 'linear', 'model__estimator__class_weight': None, 'model__estimator__C': 0.001}
'@
$ws.Range("C2").Value = $C2
$ws.Range("D2").Value = 0.1666666666666667
$ws.Range("E2").Value = '[1 1 0 0 1 0 0 0 0 1 0 1]'
$ws.Range("F2").Value = '[0 0 1 1 0 1 1 0 1 0 1 1]'
$ws.Range("G2").Value = 77
$ws.Range("H2").Value = 0.6756907701352146
$ws.Range("I2").Value = 0.04240077551254682
$ws.Range("J2").Value = 0.5602586713697824
$ws.Range("K2").Value = 0.06632229225328388

# Row 3
$A3 = @'
Pipeline(steps=[('scaler', None),
                ('selector',
                 <__main__.NamedFeatureSelector object at 0x7f3a5c78c2e0>),
                ('model',
                 BaggingClassifier(estimator=SVC(C=5, kernel='linear',
                                                 random_state=42),
                                   random_state=42))])
'@
$ws.Range("A3").Value = $A3

$ws.Range("B3").Value = 0.6095238095238096
$C3 = @'
{'selector': <__main__.NamedFeatureSelector object at 0x7f3a643ad2b0>, 'scaler': None, 'model__n_estimators': 10, 'model__estimator__kernel': 'linear', 'model__estimator__class_weight': None, 'model__estimator__C': 5}
'@
$ws.Range("C3").Value = $C3
$ws.Range("D3").Value = 0.7777777777777778
$ws.Range("E3").Value = '[1 1 0 1 0 0 1 0 1 1 1 0]'
$ws.Range("F3").Value = '[1 1 1 1 1 0 1 1 1 1 1 1]'
$ws.Range("G3").Value = 69
$ws.Range("H3").Value = 0.6746138996138996
$ws.Range("I3").Value = 0.03508057132303034
$ws.Range("J3").Value = 0.543114543114543
$ws.Range("K3").Value = 0.07128255763077197

# Row 4
$A4 = @'
Pipeline(steps=[('scaler', None),
                ('selector',
                 <__main__.NamedFeatureSelector object at 0x7f3a642efdf0>),
                ('model',
                 BaggingClassifier(estimator=SVC(C=1, class_weight='balanced',
                                                 random_state=42),
                                   n_estimators=50, random_state=42))])
'@
$ws.Range("A4").Value = $A4

$ws.Range("B4").Value = 0.6190476190476191
$C4 = @'
{'selector': <__main__.NamedFeatureSelector object at 0x7f3a64212a00>, 'scaler': None, 'model__n_estimators': 50, 'model__estimator__kernel': 'rbf', 'model__estimator__class_weight': 'balanced', 'model__estimator__C': 1}
'@
$ws.Range("C4").Value = $C4
$ws.Range("D4").Value = 0.875
$ws.Range("E4").Value = '[1 0 1 1 1 1 0 1 0 1 0 1]'
$ws.Range("F4").Value = '[1 1 1 1 1 1 0 0 0 1 0 1]'
$ws.Range("G4").Value = 42
$ws.Range("H4").Value = 0.6403880070546738
$ws.Range("I4").Value = 0.03841910271617143
$ws.Range("J4").Value = 0.5291005291005291
$ws.Range("K4").Value = 0.06539195137793538

# Row 5
$A5 = @'
Pipeline(steps=[('scaler', MinMaxScaler()),
                ('selector',
                 SelectFromModel(estimator=LinearSVC(dual=False, penalty='l1',
                                                     random_state=42))),
                ('model',
                 BaggingClassifier(estimator=SVC(C=0.0001, random_state=42),
                                   n_estimators=5, random_state=42))])
'@
$ws.Range("A5").Value = $A5

$ws.Range("B5").Value = 0.5904761904761905
$C5 = @'
{'selector': SelectFromModel(estimator=LinearSVC(dual=False, penalty='l1', random_state=42)), 'scaler': MinMaxScaler(), 'model__n_estimators': 5, 'model__estimator__kernel': 'rbf', 'model__estimator__class_weight': None, 'model__estimator__C': 0.0001}
'@
$ws.Range("C5").Value = $C5
$ws.Range("D5").Value = 0.7368421052631579
$ws.Range("E5").Value = '[1 1 0 0 0 0 1 0 1 1 1 1]'
$ws.Range("F5").Value = '[1 1 1 1 1 1 1 1 1 1 1 1]'
$ws.Range("G5").Value = 11
$ws.Range("H5").Value = 0.6385714285714286
$ws.Range("I5").Value = 0.02900492805779044
$ws.Range("J5").Value = 0.5173809523809523
$ws.Range("K5").Value = 0.05856588601692412

# Row 6
$A6 = @'
Pipeline(steps=[('scaler', None),
                ('selector',
                 SelectFromModel(estimator=ExtraTreesClassifier(random_state=42))),
                ('model',
                 BaggingClassifier(estimator=SVC(C=5, class_weight='balanced',
                                                 kernel='poly',
                                                 random_state=42),
                                   n_estimators=50, random_state=42))])
'@
$ws.Range("A6").Value = $A6

$ws.Range("B6").Value = 0.6476190476190476
$C6 = @'
{'selector': SelectFromModel(estimator=ExtraTreesClassifier(random_state=42)), 'scaler': None, 'model__n_estimators': 50, 'model__estimator__kernel': 'poly', 'model__estimator__class_weight': 'balanced', 'model__estimator__C': 5}
'@
$ws.Range("C6").Value = $C6
$ws.Range("D6").Value = 0.4615384615384615
$ws.Range("E6").Value = '[1 1 1 1 0 0 0 0 1 1 0 0]'
$ws.Range("F6").Value = '[1 1 0 0 0 1 1 1 1 0 1 0]'
$ws.Range("G6").Value = 14
$ws.Range("H6").Value = 0.6622605363984675
$ws.Range("I6").Value = 0.03657517459822247
$ws.Range("J6").Value = 0.5308155446086481
$ws.Range("K6").Value = 0.06431040025763396
